$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the "Sheet1" worksheet (with the lookup table) entirely
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete() | Out-Null

$ws = $wb.Worksheets.Item("master-template_type")

# Append the new master-data rows (92-121)
$data = @(
    ,@("RPR_UIN_CARD_TEMPLATE", "UIN card template", "eng", "superadmin", "now()")
    ,@("RPR_UIN_CARD_TEMPLATE", "قالب بطاقة UIN", "ara", "superadmin", "now()")
    ,@("RPR_UIN_CARD_TEMPLATE", "Modèle de carte UIN", "fra", "superadmin", "now()")
    ,@("RPR_UIN_DEAC_SMS", "Template for UIN Deactivation SMS", "eng", "superadmin", "now()")
    ,@("RPR_UIN_DEAC_SMS", "قالب لتعطيل UIN SMS", "ara", "superadmin", "now()")
    ,@("RPR_UIN_DEAC_SMS", "Modèle pour SMS de désactivation UIN", "fra", "superadmin", "now()")
    ,@("RPR_UIN_DEAC_EMAIL", "Template for UIN Deactivation Email", "eng", "superadmin", "now()")
    ,@("RPR_UIN_DEAC_EMAIL", "قالب لإلغاء تنشيط البريد", "ara", "superadmin", "now()")
    ,@("RPR_UIN_DEAC_EMAIL", "Modèle pour Email de désactivation UIN", "fra", "superadmin", "now()")
    ,@("RPR_UIN_REAC_SMS", "Template for UIN Reactivate SMS", "eng", "superadmin", "now()")
    ,@("RPR_UIN_REAC_SMS", "قالب لـ UIN تنشيط SMS", "ara", "superadmin", "now()")
    ,@("RPR_UIN_REAC_SMS", "Modèle pour UIN Réactiver SMS", "fra", "superadmin", "now()")
    ,@("RPR_UIN_REAC_EMAIL", "Template for UIN Reactivate Email", "eng", "superadmin", "now()")
    ,@("RPR_UIN_REAC_EMAIL", "قالب لـ UIN تنشيط البريد", "ara", "superadmin", "now()")
    ,@("RPR_UIN_REAC_EMAIL", "Modèle pour UIN Réactiver Email", "fra", "superadmin", "now()")
    ,@("reg-sms-notification", "Registration Acknowledgement Template", "eng", "superadmin", "now()")
    ,@("reg-sms-notification", "نموذج شكر التسجيل", "ara", "superadmin", "now()")
    ,@("reg-sms-notification", "accusé de réception", "fra", "superadmin", "now()")
    ,@("reg-email-notification", "Registration Acknowledgement Template", "eng", "superadmin", "now()")
    ,@("reg-email-notification", "نموذج شكر التسجيل", "ara", "superadmin", "now()")
    ,@("reg-email-notification", "accusé de réception", "fra", "superadmin", "now()")
    ,@("reg-ack-template-part1", "Registration Acknowledgement Template - Part 1", "eng", "superadmin", "now()")
    ,@("reg-ack-template-part2", "نموذج شكر التسجيل", "ara", "superadmin", "now()")
    ,@("reg-ack-template-part3", "accusé de réception", "fra", "superadmin", "now()")
    ,@("reg-ack-template-part2", "Registration Acknowledgement Template - Part 2", "eng", "superadmin", "now()")
    ,@("reg-ack-template-part3", "نموذج شكر التسجيل", "ara", "superadmin", "now()")
    ,@("reg-ack-template-part4", "accusé de réception", "fra", "superadmin", "now()")
    ,@("reg-ack-template-part3", "Registration Acknowledgement Template - Part 3", "eng", "superadmin", "now()")
    ,@("reg-ack-template-part4", "نموذج شكر التسجيل", "ara", "superadmin", "now()")
    ,@("reg-ack-template-part5", "accusé de réception", "fra", "superadmin", "now()")
)

$startRow = 92
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $r = $data[$i]
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
}

# Reflect the post-edit selection (everything below the last data row)
$lastRow = $startRow + $data.Count
$ws.Range("A" + $lastRow + ":XFD1048576").Select() | Out-Null
